{"js": "// The commit replaces the narrative paragraph that begins \"This visualization\n// uses a martini glass...\" with an updated explanation (three parameters /\n// author-driven vs user-driven content), and indents the blank paragraph\n// that follows it to match the surrounding body-text indent (0.5in / 36pt).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  paragraphs.items[i].load(\"text\");\n}\nawait context.sync();\n\n// Find the target paragraph by its distinctive leading text instead of a\n// hard-coded index, so the script is resilient to minor structural drift.\nlet targetIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"This visualization uses a martini glass\") === 0) {\n    targetIndex = i;\n    break;\n  }\n}\n\nif (targetIndex === -1) {\n  throw new Error(\"Could not locate the 'martini glass' narrative paragraph.\");\n}\n\nconst targetParagraph = paragraphs.items[targetIndex];\n\nconst newText =\n  \"This visualization uses a martini glass. \" +\n  \"There are three parameters (i.e., the three candidates I collected data from). \" +\n  \"The author-driven contents are three scenes, which user can navigate to by clicking the buttons, where the same chart type is implemented. \" +\n  \"The chart is configured based on the parameter \\u2013 the parameter decides which candidate\\u2019s data to render to the graphs. \" +\n  \"For the user driven contents, user can use the date dropdown to filter for the data of interest \\u2013 only the data in the date range are included in the graph.\";\n\n// Replace the whole paragraph's text in one shot; this keeps the run\n// formatting (Source Sans Pro / #1F1F1F) that was already on the paragraph.\ntargetParagraph.insertText(newText, Word.InsertLocation.replace);\n\n// The paragraph immediately following it (an empty spacer paragraph) gains a\n// 720-twip (0.5in / 36pt) left indent, matching the indentation used\n// elsewhere in this section.\nconst nextParagraph = paragraphs.items[targetIndex + 1];\nnextParagraph.leftIndent = 36;\n\nawait context.sync();\n", "ps1": "# The commit replaces the narrative paragraph that begins \"This visualization\n# uses a martini glass...\" with an updated explanation (three parameters /\n# author-driven vs user-driven content), and indents the blank paragraph\n# that follows it to match the surrounding body-text indent (0.5in / 36pt).\n\n$d = $word.ActiveDocument\n\n# Locate the target paragraph by its distinctive leading text instead of a\n# hard-coded index, so the script is resilient to minor structural drift.\n$searchRange = $d.Content\n$found = $searchRange.Find.Execute(\"This visualization uses a martini glass\")\nif (-not $found) {\n    throw \"Could not locate the 'martini glass' narrative paragraph.\"\n}\n\n$targetParagraph = $searchRange.Paragraphs.Item(1)\n$targetStart = $targetParagraph.Range.Start\n\n$allParagraphs = $d.Paragraphs\n$targetIndex = 0\nfor ($i = 1; $i -le $allParagraphs.Count; $i++) {\n    if ($allParagraphs.Item($i).Range.Start -eq $targetStart) {\n        $targetIndex = $i\n        break\n    }\n}\nif ($targetIndex -eq 0) {\n    throw \"Could not resolve paragraph index for the 'martini glass' paragraph.\"\n}\n\n$newText = \"This visualization uses a martini glass. \" + `\n    \"There are three parameters (i.e., the three candidates I collected data from). \" + `\n    \"The author-driven contents are three scenes, which user can navigate to by clicking the buttons, where the same chart type is implemented. \" + `\n    \"The chart is configured based on the parameter \" + [char]0x2013 + \" the parameter decides which candidate\" + [char]0x2019 + \"s data to render to the graphs. \" + `\n    \"For the user driven contents, user can use the date dropdown to filter for the data of interest \" + [char]0x2013 + \" only the data in the date range are included in the graph.\"\n\n# Replace the whole paragraph's text (without its trailing paragraph mark) in\n# one shot; this keeps the run formatting (Source Sans Pro / #1F1F1F) that was\n# already on the paragraph. Re-wrap the paragraph's Range through Document.Range\n# so the assignment replaces the full span instead of just its first word.\n$pRange = $allParagraphs.Item($targetIndex).Range\n$fullRange = $d.Range($pRange.Start, $pRange.End)\n$fullRange.Text = $newText\n\n# The paragraph immediately following it (an empty spacer paragraph) gains a\n# 720-twip (0.5in / 36pt) left indent, matching the indentation used\n# elsewhere in this section.\n$nextParagraph = $d.Paragraphs.Item($targetIndex + 1)\n$nextParagraph.Format.LeftIndent = 36\n"}
